$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = "c76eb7bc1a2e6e67034fc7d750762de9"
$ws.Range("B17").Value = "6d8ffd1d66c53c710be8772851e1d28a"
$ws.Range("B24").Value = "f922ed9e447644263a1a97de707e5cb8"
$ws.Range("B136").Value = "b7039b57dbda92005e340256ad999b90"
$ws.Range("B159").Value = "e180276079263c04640119ac2f9a2356"
$ws.Range("B169").Value = "57c8ebb0b1bfe05484cfbeee6e543676"
$ws.Range("B183").Value = "1566ad624c9b683444f8640e7090cefd"
$ws.Range("B200").Value = "d5ef55e36803ff9c65c83cdd13fffe52"
$ws.Range("B228").Value = "fe38701a3da4b84079059572acfcc9b3"
$ws.Range("B246").Value = "071af2aceba94931a20530f3df305f55"
$ws.Range("B276").Value = "6ce9b456d4485a7c5da99bf32a286582"
$ws.Range("B339").Value = "099ff95134ac2a6dda1c6112387b1c53"
$ws.Range("B411").Value = "fbd76eead3baf2a33ed1e1dab1fb7e73"
$ws.Range("B420").Value = "bf3569543f5afe0bd329968445d710df"
$ws.Range("B448").Value = "dbd283247f7d76505c7c28daa379885d"
$ws.Range("B464").Value = "88ca15026fa327f90edcf2607339c165"
$ws.Range("B507").Value = "28b7f4082aa807fa960d3091d6953006"
$ws.Range("B508").Value = "3bb24bf20af84bd73d4fd48e30da03f3"
$ws.Range("B521").Value = "3962d32114f3fb69ae6f12f86a119019"
$ws.Range("B522").Value = "a36d455db57eed089334afd0313fbd4c"
$ws.Range("B532").Value = "320c9d5b1e38d46bf285d4beb72f820c"
$ws.Range("B555").Value = "94c8a699ba72fa2ba49483e62eaeeb5b"
$ws.Range("B574").Value = "c2773ef09b571a4d55e3f514b1138e7d"
$ws.Range("B580").Value = "521ce29e8304ca26acab34907e3d08da"
$ws.Range("B624").Value = "8eed330081db7ea415c2ac50c2458014"
$ws.Range("B635").Value = "d450c3da6f90944d2dbd85eeeee6c17e"
$ws.Range("B673").Value = "2ede366eee4394e48ea0925f9464345c"
$ws.Range("B674").Value = "654c1ba0472b17af82efd250300ae113"
$ws.Range("B708").Value = "12e5dbeb119384264be0298d3ffb04dd"
$ws.Range("B712").Value = "c3305368066951b035b3eec49bbfc9ce"
$ws.Range("B723").Value = "3d55dde6eea0e77c61e852a4347905de"
$ws.Range("B734").Value = "06fceded922be85cd11a7b87e4c5a2c0"
$ws.Range("B750").Value = "bebe597650251d7dc4b5abfc624cebb2"
$ws.Range("B764").Value = "d779807ca271b03402b356011c198692"
$ws.Range("B769").Value = "753333d5fc4f3ef466f64ee800da8620"
$ws.Range("B794").Value = "1bdfab8e7202f5daefeedaa98f3e8aef"
$ws.Range("B827").Value = "6f14a86add7ba4c658e6672d743c2b75"
$ws.Range("B838").Value = "10e0d3fcba82c94ccc94802d6c5c9179"
$ws.Range("B843").Value = "08ec81e9257330f99b6ec686fc7b6d56"
$ws.Range("B863").Value = "c61c485da4221da22910550d738db2b2"
$ws.Range("B877").Value = "9ec6e776bd4b0df7de15559b3bc14cdb"
$ws.Range("B913").Value = "d3525cdd28c8b00d6338a12f4b459ddd"
$ws.Range("B937").Value = "0f0e708f250eb6be44ce4686d1174aa4"
